{"js": "// Apply the dated worksheet's text updates: the header date and the\n// twenty-five \"N\u00d7N=N\" multiplication answers scattered throughout the\n// table cells. Each old string is unique in the document, so a scoped\n// case-sensitive search-and-replace for every (old, new) pair reproduces\n// the diff exactly without disturbing formatting/runs.\nconst replacements = [\n  [\"2025-10-10 Friday\", \"2025-10-11 Saturday\"],\n  [\"887\u00d78=7096\", \"201\u00d76=1206\"],\n  [\"953\u00d75=4765\", \"257\u00d76=1542\"],\n  [\"463\u00d74=1852\", \"411\u00d77=2877\"],\n  [\"141\u00d74=564\", \"191\u00d72=382\"],\n  [\"579\u00d75=2895\", \"580\u00d72=1160\"],\n  [\"169\u00d75=845\", \"643\u00d76=3858\"],\n  [\"965\u00d77=6755\", \"333\u00d79=2997\"],\n  [\"407\u00d74=1628\", \"771\u00d75=3855\"],\n  [\"195\u00d78=1560\", \"466\u00d79=4194\"],\n  [\"186\u00d74=744\", \"551\u00d73=1653\"],\n  [\"807\u00d77=5649\", \"503\u00d74=2012\"],\n  [\"479\u00d79=4311\", \"282\u00d76=1692\"],\n  [\"566\u00d78=4528\", \"863\u00d77=6041\"],\n  [\"857\u00d74=3428\", \"371\u00d75=1855\"],\n  [\"716\u00d78=5728\", \"684\u00d73=2052\"],\n  [\"794\u00d75=3970\", \"690\u00d75=3450\"],\n  [\"906\u00d78=7248\", \"841\u00d76=5046\"],\n  [\"887\u00d75=4435\", \"482\u00d79=4338\"],\n  [\"616\u00d72=1232\", \"170\u00d72=340\"],\n  [\"868\u00d79=7812\", \"192\u00d76=1152\"],\n  [\"944\u00d75=4720\", \"524\u00d72=1048\"],\n  [\"632\u00d79=5688\", \"119\u00d79=1071\"],\n  [\"133\u00d77=931\", \"507\u00d74=2028\"],\n  [\"859\u00d72=1718\", \"316\u00d74=1264\"],\n  [\"407\u00d75=2035\", \"316\u00d79=2844\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the dated worksheet's text updates: the header date and the\n# twenty-five \"N\u00d7N=N\" multiplication answers scattered throughout the\n# table cells. Each old string is unique in the document, so a\n# Find/Replace-All pass for every (old, new) pair over the whole document\n# content reproduces the diff exactly without disturbing formatting/runs.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2025-10-10 Friday\", \"2025-10-11 Saturday\"),\n  @(\"887\u00d78=7096\", \"201\u00d76=1206\"),\n  @(\"953\u00d75=4765\", \"257\u00d76=1542\"),\n  @(\"463\u00d74=1852\", \"411\u00d77=2877\"),\n  @(\"141\u00d74=564\", \"191\u00d72=382\"),\n  @(\"579\u00d75=2895\", \"580\u00d72=1160\"),\n  @(\"169\u00d75=845\", \"643\u00d76=3858\"),\n  @(\"965\u00d77=6755\", \"333\u00d79=2997\"),\n  @(\"407\u00d74=1628\", \"771\u00d75=3855\"),\n  @(\"195\u00d78=1560\", \"466\u00d79=4194\"),\n  @(\"186\u00d74=744\", \"551\u00d73=1653\"),\n  @(\"807\u00d77=5649\", \"503\u00d74=2012\"),\n  @(\"479\u00d79=4311\", \"282\u00d76=1692\"),\n  @(\"566\u00d78=4528\", \"863\u00d77=6041\"),\n  @(\"857\u00d74=3428\", \"371\u00d75=1855\"),\n  @(\"716\u00d78=5728\", \"684\u00d73=2052\"),\n  @(\"794\u00d75=3970\", \"690\u00d75=3450\"),\n  @(\"906\u00d78=7248\", \"841\u00d76=5046\"),\n  @(\"887\u00d75=4435\", \"482\u00d79=4338\"),\n  @(\"616\u00d72=1232\", \"170\u00d72=340\"),\n  @(\"868\u00d79=7812\", \"192\u00d76=1152\"),\n  @(\"944\u00d75=4720\", \"524\u00d72=1048\"),\n  @(\"632\u00d79=5688\", \"119\u00d79=1071\"),\n  @(\"133\u00d77=931\", \"507\u00d74=2028\"),\n  @(\"859\u00d72=1718\", \"316\u00d74=1264\"),\n  @(\"407\u00d75=2035\", \"316\u00d79=2844\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
